# TC_RunManager.xlsx - Batch 5 scripts: Namemapping and Datasheet
# Adds three new rows (58-60) to the GlobalTestCase sheet describing the
# new "ReverseGL" / "CopyGL" / "CurrencyGL" units, and updates the sheet's
# view (scroll position / active selection) to match where the author was
# working when the rows were added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalTestCase")
$ws.Activate()

# --- New data rows -------------------------------------------------------
# Row 55 (style 3 in cols A-D, style 2 in E-I) is the closest existing
# pattern; build each new row cell-by-cell so every column picks up the
# exact formatting the workbook already uses (A-C + E = "style 3",
# D + F-I = "style 2"), without minting any new style entries.

function Set-RowFormatting($rowNum) {
    $ws.Range("A55:C55").Copy() | Out-Null
    $ws.Range("A$rowNum`:C$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("D2").Copy() | Out-Null
    $ws.Range("D$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("A55").Copy() | Out-Null
    $ws.Range("E$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("D2").Copy() | Out-Null
    $ws.Range("F$rowNum`:I$rowNum").PasteSpecial(-4122) | Out-Null
}

Set-RowFormatting 58
Set-RowFormatting 59
Set-RowFormatting 60

$excel.CutCopyMode = $false

# Row 58 - ReverseGL / Reverse a general ledger
$ws.Range("A58").Value = "ReverseGL"
$ws.Range("B58").Value = "ReverseGL"
$ws.Range("C58").Value = "Reverse a general ledger"
$ws.Range("D58").Value = "No"
$ws.Range("E58").Value = "No"
$ws.Range("F58").Value = "No"
$ws.Range("G58").Value = "No"
$ws.Range("H58").Value = "No"
$ws.Range("I58").Value = "No"

# Row 59 - CopyGL / Copy a general Ledger
$ws.Range("A59").Value = "CopyGL"
$ws.Range("B59").Value = "CopyGL"
$ws.Range("C59").Value = "Copy a general Ledger"
$ws.Range("D59").Value = "No"
$ws.Range("E59").Value = "No"
$ws.Range("F59").Value = "No"
$ws.Range("G59").Value = "No"
$ws.Range("H59").Value = "No"
$ws.Range("I59").Value = "No"

# Row 60 - CurrencyGL / CreateCurrencyJournal / Create a Currency Journal
$ws.Range("A60").Value = "CurrencyGL"
$ws.Range("B60").Value = "CreateCurrencyJournal"
$ws.Range("C60").Value = "Create a Currency Journal"
$ws.Range("D60").Value = "No"
$ws.Range("E60").Value = "Yes"
$ws.Range("F60").Value = "No"
$ws.Range("G60").Value = "No"
$ws.Range("H60").Value = "No"
$ws.Range("I60").Value = "No"

# --- View state ------------------------------------------------------------
# Author scrolled down to row 43 and left the active cell on C62.
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C62").Select() | Out-Null
